# Swap the species/record data between row 2 and row 3 for the columns
# that differ: A, B, E, F, G, H, Q, R. All other columns are identical
# between the two rows already, so no other changes are required.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell2 = $ws.Range("$col" + "2")
    $cell3 = $ws.Range("$col" + "3")

    $val2 = $cell2.Value()
    $val3 = $cell3.Value()

    $cell2.Value = $val3
    $cell3.Value = $val2
}
